$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 3-27 (only cells whose value actually changes) ---
$ws.Range("B3").Value = 0.7941176470588235
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 27
$ws.Range("H3").Value = 7
$ws.Range("J3").Value = "best"
$ws.Range("K3").Value = 0.9491525423728814
$ws.Range("L3").Value = 56
$ws.Range("M3").Value = 56
$ws.Range("Q3").Value = 3
$ws.Range("B4").Value = 0.7222222222222222
$ws.Range("C4").Value = 26
$ws.Range("D4").Value = 26
$ws.Range("H4").Value = 10
$ws.Range("J4").Value = "interesting"
$ws.Range("K4").Value = 0.9090909090909091
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = 30
$ws.Range("Q4").Value = 3
$ws.Range("B5").Value = 0.5993150684931506
$ws.Range("C5").Value = 175
$ws.Range("D5").Value = 175
$ws.Range("H5").Value = 117
$ws.Range("K5").Value = 0.8913043478260869
$ws.Range("L5").Value = 41
$ws.Range("M5").Value = 41
$ws.Range("Q5").Value = 5
$ws.Range("B6").Value = 0.2151162790697674
$ws.Range("C6").Value = 111
$ws.Range("D6").Value = 111
$ws.Range("H6").Value = 405
$ws.Range("K6").Value = 0.8392857142857143
$ws.Range("L6").Value = 94
$ws.Range("M6").Value = 94
$ws.Range("Q6").Value = 18
$ws.Range("B7").Value = 0.1798941798941799
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 34
$ws.Range("H7").Value = 155
$ws.Range("J7").Value = "thanks"
$ws.Range("K7").Value = 0.8292682926829268
$ws.Range("L7").Value = 68
$ws.Range("M7").Value = 68
$ws.Range("Q7").Value = 14
$ws.Range("B8").Value = 0.0992063492063492
$ws.Range("C8").Value = 25
$ws.Range("D8").Value = 25
$ws.Range("H8").Value = 227
$ws.Range("J8").Value = "positive"
$ws.Range("K8").Value = 0.8103448275862069
$ws.Range("Q8").Value = 11
$ws.Range("J9").Value = "free"
$ws.Range("K9").Value = 0.8
$ws.Range("L9").Value = 96
$ws.Range("M9").Value = 96
$ws.Range("Q9").Value = 24
$ws.Range("J10").Value = "special"
$ws.Range("K10").Value = 0.7777777777777778
$ws.Range("L10").Value = 28
$ws.Range("M10").Value = 28
$ws.Range("Q10").Value = 8
$ws.Range("J11").Value = "thank"
$ws.Range("K11").Value = 0.7421875
$ws.Range("L11").Value = 95
$ws.Range("M11").Value = 95
$ws.Range("Q11").Value = 33
$ws.Range("J12").Value = "safe"
$ws.Range("K12").Value = 0.7183098591549296
$ws.Range("L12").Value = 102
$ws.Range("M12").Value = 102
$ws.Range("Q12").Value = 40
$ws.Range("J13").Value = "support"
$ws.Range("K13").Value = 0.7075471698113207
$ws.Range("L13").Value = 75
$ws.Range("M13").Value = 75
$ws.Range("Q13").Value = 31
$ws.Range("J14").Value = "safety"
$ws.Range("K14").Value = 0.7058823529411765
$ws.Range("L14").Value = 36
$ws.Range("M14").Value = 36
$ws.Range("Q14").Value = 15
$ws.Range("J15").Value = "good"
$ws.Range("K15").Value = 0.65625
$ws.Range("L15").Value = 105
$ws.Range("M15").Value = 105
$ws.Range("Q15").Value = 55
$ws.Range("J16").Value = "better"
$ws.Range("K16").Value = 0.6507936507936508
$ws.Range("L16").Value = 41
$ws.Range("M16").Value = 41
$ws.Range("Q16").Value = 22
$ws.Range("J17").Value = "fresh"
$ws.Range("K17").Value = 0.6041666666666666
$ws.Range("L17").Value = 29
$ws.Range("M17").Value = 29
$ws.Range("Q17").Value = 19
$ws.Range("J18").Value = "relief"
$ws.Range("K18").Value = 0.6
$ws.Range("L18").Value = 30
$ws.Range("M18").Value = 30
$ws.Range("Q18").Value = 20
$ws.Range("J19").Value = "well"
$ws.Range("K19").Value = 0.5531914893617021
$ws.Range("L19").Value = 52
$ws.Range("M19").Value = 52
$ws.Range("Q19").Value = 42
$ws.Range("J20").Value = "hand"
$ws.Range("K20").Value = 0.5143603133159269
$ws.Range("L20").Value = 197
$ws.Range("M20").Value = 197
$ws.Range("Q20").Value = 186
$ws.Range("J21").Value = "like"
$ws.Range("K21").Value = 0.5058823529411764
$ws.Range("L21").Value = 172
$ws.Range("M21").Value = 172
$ws.Range("Q21").Value = 168
$ws.Range("J22").Value = "care"
$ws.Range("K22").Value = 0.5056179775280899
$ws.Range("L22").Value = 45
$ws.Range("M22").Value = 45
$ws.Range("Q22").Value = 44
$ws.Range("K23").Value = 0.4440677966101695
$ws.Range("L23").Value = 131
$ws.Range("M23").Value = 131
$ws.Range("Q23").Value = 164
$ws.Range("J24").Value = "protect"
$ws.Range("K24").Value = 0.410958904109589
$ws.Range("L24").Value = 30
$ws.Range("M24").Value = 30
$ws.Range("Q24").Value = 43
$ws.Range("J25").Value = "hope"
$ws.Range("K25").Value = 0.4
$ws.Range("L25").Value = 26
$ws.Range("M25").Value = 26
$ws.Range("Q25").Value = 39
$ws.Range("J26").Value = "increase"
$ws.Range("K26").Value = 0.3846153846153846
$ws.Range("L26").Value = 30
$ws.Range("M26").Value = 30
$ws.Range("Q26").Value = 48
$ws.Range("J27").Value = "please"
$ws.Range("K27").Value = 0.3723849372384937
$ws.Range("L27").Value = 89
$ws.Range("M27").Value = 89
$ws.Range("Q27").Value = 150

# --- Append new rows 28-31 ---
$ws.Range("J28").Value = "store"
$ws.Range("K28").Value = 0.05257270693512305
$ws.Range("L28").Value = 47
$ws.Range("M28").Value = 47
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 847
$ws.Range("J29").Value = "grocery"
$ws.Range("K29").Value = 0.02996670366259711
$ws.Range("L29").Value = 27
$ws.Range("M29").Value = 27
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 874
$ws.Range("J30").Formula = '=TEXT(19,"0")'
$ws.Range("J30").Copy()
$ws.Range("J30").PasteSpecial(-4163)
$ws.Range("K30").Value = 0.01308411214953271
$ws.Range("L30").Value = 28
$ws.Range("M30").Value = 31
$ws.Range("N30").Value = 0.9
$ws.Range("O30").Value = 0.09999999999999998
$ws.Range("P30").Value = $true
$ws.Range("Q30").Value = 2112
$ws.Range("J31").Value = "co"
$ws.Range("K31").Value = 0.01128668171557562
$ws.Range("L31").Value = 35
$ws.Range("M31").Value = 41
$ws.Range("N31").Value = 0.85
$ws.Range("O31").Value = 0.15
$ws.Range("P31").Value = $true
$ws.Range("Q31").Value = 3066

# --- Copy J-column style (bold + border + centered) from row 27 onto new J28:J31 ---
$ws.Range("J27").Copy()
$ws.Range("J28:J31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
